$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 -> RECURRENTE_100K-200K : append new client id 60299
$ws.Range("B6").Value = $ws.Range("B6").Value2 + ".60299"

# Row 5 -> RECURRENTE_<_100K : append new client id 60245
$ws.Range("B5").Value = $ws.Range("B5").Value2 + ".60245"
